# Scheduled market-data refresh: update crafting-leve profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) across all
# eight job sheets with freshly scraped prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H8").Value = 33.076923
$ws.Range("I8").Value = 36.363636
$ws.Range("K8").Value = 109.090908
$ws.Range("M8").Value = 29.909092

$ws.Range("H16").Value = 17003.334
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 23005
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 23005
$ws.Range("M16").Value = -4770
$ws.Range("N16").Value = -23465

$ws.Range("H40").Value = 1965.6666
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H76").Value = 4718.2856
$ws.Range("I76").Value = 4187.3335
$ws.Range("K76").Value = 4187.3335
$ws.Range("M76").Value = -3872.3335

$ws.Range("H79").Value = 4718.2856
$ws.Range("I79").Value = 4187.3335
$ws.Range("K79").Value = 4187.3335
$ws.Range("M79").Value = -3095.3335

$ws.Range("H132").Value = 2288.1
$ws.Range("I132").Value = 2288.1
$ws.Range("K132").Value = 6864.299999999999
$ws.Range("M132").Value = -4334.299999999999

$ws = $wb.Worksheets("ARM")
$ws.Range("H63").Value = 3927.5715
$ws.Range("I63").Value = 3623.25
$ws.Range("K63").Value = 3623.25
$ws.Range("M63").Value = -2937.25

$ws.Range("H66").Value = 3927.5715
$ws.Range("I66").Value = 3623.25
$ws.Range("K66").Value = 18116.25
$ws.Range("M66").Value = -14684.25

$ws.Range("H98").Value = 12450.833
$ws.Range("J98").Value = 12450.833
$ws.Range("L98").Value = 12450.833
$ws.Range("N98").Value = -18440.833

$ws = $wb.Worksheets("BSM")
$ws.Range("H105").Value = 4127.5713
$ws.Range("I105").Value = 3982.3333
$ws.Range("K105").Value = 3982.3333
$ws.Range("M105").Value = -2235.3333

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 388.33334
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 416
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 416
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -1116

$ws.Range("H31").Value = 2137.0588
$ws.Range("I31").Value = 1645.625
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 1645.625
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -1350.625
$ws.Range("N31").Value = -10590

$ws.Range("H34").Value = 2137.0588
$ws.Range("I34").Value = 1645.625
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 1645.625
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -1443.625
$ws.Range("N34").Value = -10404

$ws.Range("H107").Value = 739.6667
$ws.Range("I107").Value = 224
$ws.Range("K107").Value = 224
$ws.Range("M107").Value = 1696

$ws.Range("H134").Value = 2462.7273
$ws.Range("I134").Value = 2499
$ws.Range("K134").Value = 7497
$ws.Range("M134").Value = -4962

$ws = $wb.Worksheets("CUL")
$ws.Range("H9").Value = 557.1429000000001
$ws.Range("J9").Value = 557.1429000000001
$ws.Range("L9").Value = 1671.4287
$ws.Range("N9").Value = -2119.4287

$ws.Range("H44").Value = 1376.6
$ws.Range("I44").Value = 1545.75
$ws.Range("J44").Value = 700
$ws.Range("K44").Value = 4637.25
$ws.Range("L44").Value = 2100
$ws.Range("M44").Value = -4239.25
$ws.Range("N44").Value = -2896

$ws.Range("H51").Value = 999.75
$ws.Range("I51").Value = 994.5
$ws.Range("K51").Value = 2983.5
$ws.Range("M51").Value = -2523.5

$ws.Range("H55").Value = 1880
$ws.Range("I55").Value = 600
$ws.Range("K55").Value = 1800
$ws.Range("M55").Value = -1623

$ws.Range("H60").Value = 2000
$ws.Range("J60").Value = 2000
$ws.Range("L60").Value = 6000
$ws.Range("N60").Value = -6502

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws = $wb.Worksheets("GSM")
$ws.Range("H12").Value = 7334.6665
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 10502
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 10502
$ws.Range("M12").Value = -860
$ws.Range("N12").Value = -10782

$ws.Range("H48").Value = 16245
$ws.Range("J48").Value = 17500
$ws.Range("L48").Value = 17500
$ws.Range("N48").Value = -18470

$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

$ws.Range("H80").Value = 2156.8
$ws.Range("I80").Value = 1971.25
$ws.Range("J80").Value = 2899
$ws.Range("K80").Value = 1971.25
$ws.Range("L80").Value = 2899
$ws.Range("M80").Value = -973.25
$ws.Range("N80").Value = -4895

$ws.Range("H83").Value = 2156.8
$ws.Range("I83").Value = 1971.25
$ws.Range("J83").Value = 2899
$ws.Range("K83").Value = 9856.25
$ws.Range("L83").Value = 14495
$ws.Range("M83").Value = -4864.25
$ws.Range("N83").Value = -24479

$ws.Range("H107").Value = 237
$ws.Range("I107").Value = 237
$ws.Range("K107").Value = 237
$ws.Range("M107").Value = 1683

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 17384.5
$ws.Range("I7").Value = 18365.941
$ws.Range("J7").Value = 700
$ws.Range("K7").Value = 18365.941
$ws.Range("L7").Value = 700
$ws.Range("M7").Value = -18253.941
$ws.Range("N7").Value = -924

$ws.Range("H46").Value = 3079.56
$ws.Range("I46").Value = 1298.75
$ws.Range("J46").Value = 3418.762
$ws.Range("K46").Value = 1298.75
$ws.Range("L46").Value = 3418.762
$ws.Range("M46").Value = -1110.75
$ws.Range("N46").Value = -3794.762

$ws.Range("H105").Value = 30666.334
$ws.Range("J105").Value = 30666.334
$ws.Range("L105").Value = 30666.334
$ws.Range("N105").Value = -37654.334

$ws.Range("H126").Value = 17384.5
$ws.Range("I126").Value = 18365.941
$ws.Range("J126").Value = 700
$ws.Range("K126").Value = 55097.823
$ws.Range("L126").Value = 2100
$ws.Range("M126").Value = -52627.823
$ws.Range("N126").Value = -7040

$ws = $wb.Worksheets("WVR")
$ws.Range("H2").Value = 14001
$ws.Range("I2").Value = 8002
$ws.Range("J2").Value = 20000
$ws.Range("K2").Value = 8002
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = -7890
$ws.Range("N2").Value = -20224

$ws.Range("H41").Value = 35993.8
$ws.Range("J41").Value = 32492.25
$ws.Range("L41").Value = 32492.25
$ws.Range("N41").Value = -33272.25

$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

$ws.Range("H122").Value = 3321.9473
$ws.Range("I122").Value = 2732.3076
$ws.Range("K122").Value = 8196.9228
$ws.Range("M122").Value = -5746.9228

$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
